$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: kareem -> Not Paid, total 0, payment method cleared
$ws.Range("K2").Value = "Not Paid"
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = ""

# Row 3: nada -> Paid, total 1006, payment method Vodafone Cash
$ws.Range("K3").Value = "Paid"
$ws.Range("L3").Value = 1006
$ws.Range("N3").Value = "Vodafone Cash"

# Row 4: youisf -> total updated to 1006 (status/method unchanged)
$ws.Range("L4").Value = 1006

# Row 5: -> Not Paid, total 0, payment method cleared
$ws.Range("K5").Value = "Not Paid"
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = ""
